# fix a save bug
# Append two new log rows (57, 58) to the bottom of the time-tracking
# table on Sheet1, and move the active selection onto the new last
# entry (B58), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 57 -----------------------------------------------------------
# "2012.7.30" looks like a date to Excel's smart-entry parser, so it
# would silently become a serial-number date cell instead of staying
# text. Stage it in a scratch cell that is explicitly text-formatted,
# then move the *value* into place and restore the destination cell's
# original (General) formatting by copying it from the row above -
# this keeps A57 a plain shared-string cell like every other date
# label in the sheet.
$ws.Range("ZZ1").NumberFormat = "@"
$ws.Range("ZZ1").Value = "2012.7.30"
$ws.Range("ZZ1").Copy()
$ws.Range("A57").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("ZZ1").Clear()

$ws.Range("B57").Value = "增加较色切换  累计得分"
$ws.Range("B56").Copy()
$ws.Range("B57").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D57").Value = 3

# --- Row 58 -------------------------------------------------------------
$ws.Range("A58").Value = "2012.7.31+8.1"
$ws.Range("A56").Copy()
$ws.Range("A58").PasteSpecial(-4122)

$ws.Range("B58").Value = "增加商店"
$ws.Range("B56").Copy()
$ws.Range("B58").PasteSpecial(-4122)

$ws.Range("C58").Value = "存储数据老被清0.或读错。是不是没有初始化？"
$ws.Range("B56").Copy()
$ws.Range("C58").PasteSpecial(-4122)

$ws.Range("D58").Value = 10

# --- Selection ----------------------------------------------------------
$ws.Range("B58").Select()
